# Update the "Gaz" sheet: add Last Price / Last Volume / End of Day Index columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gaz")

# New header cells (row 1), styled like the existing headers (bold, centered, bordered).
$ws.Range("E1").Value = "Last Price"
$ws.Range("F1").Value = "Last Volume"
$ws.Range("G1").Value = "End of Day Index"

$ws.Range("E1:G1").Font.Bold = $true
$ws.Range("E1:G1").HorizontalAlignment = -4108
$ws.Range("E1:G1").VerticalAlignment = -4160
$ws.Range("E1:G1").Borders.LineStyle = 1

# Row 2 gets blank placeholder cells in the new columns (materialize the cells
# without giving them a value or a lingering style).
$ws.Range("E2:G2").Borders.LineStyle = 1
$ws.Range("E2:G2").Borders.LineStyle = -4142

# Row 3 carries the actual data values.
$ws.Range("E3").Value = 38.95
$ws.Range("F3").Value = 24000
$ws.Range("G3").Value = 38.201
